# Entsoe "Actual Production Hydro Water Reservoir" daily refresh:
# shift every timestamp in column A forward by one day, and replace the
# "Actual Production (MW)" values in column B with the newly fetched
# readings for the next day (row 2 maps to the new oldest timestamp,
# row 193 to the newest).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newB = @(333,329,328,326,328,329,330,328,0,327,328,326,327,0,326,327,326,325,321,311,336,340,0,0,483,433,0,436,580,518,509,512,375,374,447,449,293,279,276,274,299,308,306,309,216,208,209,207,227,231,234,215,319,298,294,336,494,512,526,536,617,603,602,0,637,648,664,663,686,687,684,683,650,767,762,769,648,0,731,0,682,670,671,669,599,594,609,598,358,349,351,356,380,379,377,381,378,377,0,368,358,357,0,359,362,359,360,359,0,361,0,360,363,361,358,360,361,360,359,365,396,395,396,397,407,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $newB.Length; $i++) {
    $r = $i + 2
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Value = $cellA.Value2 + 1
    $ws.Cells.Item($r, 2).Value = $newB[$i]
}
